# The deck ships with two theme parts whose roles are "swapped" relative to
# a normal PowerPoint save: ppt/theme/theme1.xml (the "Office Theme" colours)
# is wired to the Notes Master, while ppt/theme/theme2.xml (the "Integral" /
# "Red Violet" colours) is wired to the Slide Master that the slides actually
# use. The edit re-saves the presentation so the Slide Master's theme becomes
# the plain "Office Theme" colour scheme (previously only used by the Notes
# Master), swapping the two themes' colour content.
#
# We drive this through the Slide Master's live ThemeColorScheme, writing
# each of the twelve theme colour slots (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) to the target "Office Theme" RGB values.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# index -> (name, target RGB as a PowerPoint BGR-packed long: R + G*256 + B*65536)
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
